$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the latest cryptos-list refresh.
# Plain numeric-looking text values need to be forced back to Text
# (NumberFormat "@") after assignment, then restyled to "Normal" so
# no stray cell style index is left behind (these cells carry no
# explicit style in the source workbook).

$ws.Range("D2").Value = '35.119.78'
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("D3").Value = '1.810.76'
$ws.Range("E3").Value = '  -2.35%  '
$ws.Range("E4").Value = '  +0.76%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.99%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("E7").Value = '  +0.72%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.34'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.58%  '
$ws.Range("E9").Value = '  +5.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0684'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.30%  '
$ws.Range("E11").Value = '  -0.60%  '
$ws.Range("D12").Value = '2.074.58'
$ws.Range("E12").Value = '  -2.06%  '
$ws.Range("D13").Value = '1.804.01'
$ws.Range("E13").Value = '  -2.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.664'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '11.02'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.64'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.13%  '
$ws.Range("D17").Value = '35.084.83'
$ws.Range("E17").Value = '  -0.92%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.63'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("D19").Value = '0.0₃0790'
$ws.Range("E19").Value = '  -0.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '238.31'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.93'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.70'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.29%  '
$ws.Range("E23").Value = '  +0.75%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '172.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.84'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.50'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.29%  '
$ws.Range("E28").Value = '  -1.71%  '
$ws.Range("E29").Value = '  +20.14%  '
$ws.Range("E30").Value = '  +0.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.13'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0554'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.02'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.92%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.77'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.15'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '91.89'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.46%  '
$ws.Range("E37").Value = '  +0.27%  '
$ws.Range("E38").Value = '  -0.80%  '
$ws.Range("D39").Value = '1.313.67'
$ws.Range("E39").Value = '  -2.21%  '
$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.28'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.11%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.86%  '
$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.47'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.88%  '
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '14.50'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.92%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.29'
$ws.Range("D44").Style = "Normal"
$ws.Range("E45").Value = '  -2.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.29'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.90%  '
$ws.Range("E47").Value = '  -1.54%  '
$ws.Range("D48").Value = '1.991.41'
$ws.Range("E48").Value = '  -1.12%  '
$ws.Range("E49").Value = '  +0.75%  '
$ws.Range("E50").Value = '  +4.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '99.26'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.60%  '
